$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Update the "last saved" date shown in the page header (SAVEDATE field
#    result) from 26.01.2023 to 02.02.2023.
# ---------------------------------------------------------------------------
$hdr = $d.Sections.Item(1).Headers.Item(1)
$hdr.Range.Find.Execute("26.01.2023", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "02.02.2023", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Add a new (still empty/un-worded) "Aufgabe" heading paragraph - this is
#    the paragraph style "Aufgabe_mit_Nr" (AufgabemitNr) which is numbered
#    automatically ("Aufgabe %1"). It is inserted right after the bicycle
#    task ("Erstelle eine Instanz der Klasse "Fahrrad" und rufe die Methoden
#    auf.") and its following blank paragraph, i.e. directly before the next
#    task about the "Zweirad" parent class. Because of the numbering, this
#    new heading will automatically show up as "Aufgabe 29".
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like '*Fahrrad" und rufe die Methoden auf.*') {
        $target = $i
        break
    }
}

# The task text paragraph is directly followed by one blank paragraph; the
# new heading paragraph needs to be inserted right after that blank one.
$blank = $d.Paragraphs.Item($target + 1)
$blank.Range.InsertParagraphAfter() | Out-Null

$newHeading = $d.Paragraphs.Item($target + 2)
$newHeading.Style = "AufgabemitNr"
$newHeading.Format.LeftIndent = 0
$newHeading.Format.FirstLineIndent = 0
